$d = $word.ActiveDocument

# 1. "knowledge on how" -> "knowledge of how"
$d.Content.Find.Execute("knowledge on how", $false, $false, $false, $false, $false, $true, 1, $false, "knowledge of how", 2) | Out-Null

# 2. "Due to such revolution" -> "Due to such a revolution"
$d.Content.Find.Execute("Due to such revolution", $false, $false, $false, $false, $false, $true, 1, $false, "Due to such a revolution", 2) | Out-Null

# 3. "main aim for this research" -> "main aim of this research"
$d.Content.Find.Execute("main aim for this research", $false, $false, $false, $false, $false, $true, 1, $false, "main aim of this research", 2) | Out-Null

# 4. "is weather a microcomputer" -> "is whether a microcomputer"
$d.Content.Find.Execute("is weather a microcomputer", $false, $false, $false, $false, $false, $true, 1, $false, "is whether a microcomputer", 2) | Out-Null

# 5. " UAV's. Such devices " -> " UAVs. Such devices " (remove apostrophe)
$d.Content.Find.Execute(" UAV" + [char]0x2019 + "s. Such devices ", $false, $false, $false, $false, $false, $true, 1, $false, " UAVs. Such devices ", 2) | Out-Null

# 6. "since, it could" -> "since it could" (remove comma)
$d.Content.Find.Execute("since, it could", $false, $false, $false, $false, $false, $true, 1, $false, "since it could", 2) | Out-Null

# 7. "automates" -> "automate"
$d.Content.Find.Execute("automates", $false, $false, $false, $false, $false, $true, 1, $false, "automate", 2) | Out-Null

# 8. "final part within this dissertation" -> "final part of this dissertation"
$d.Content.Find.Execute("final part within this dissertation", $false, $false, $false, $false, $false, $true, 1, $false, "final part of this dissertation", 2) | Out-Null

# 9. "discussion of  results" (double space) -> "discussion of results" (single space)
$d.Content.Find.Execute("discussion of  results", $false, $false, $false, $false, $false, $true, 1, $false, "discussion of results", 2) | Out-Null

# 10. "data collection will be " -> "data collection will be" (remove trailing space, leaving the lone-space run that follows)
$d.Content.Find.Execute("data collection will be ", $false, $false, $false, $false, $false, $true, 1, $false, "data collection will be", 2) | Out-Null

# 11. "compare results of different" -> "compare the results of different"
$d.Content.Find.Execute("compare results of different", $false, $false, $false, $false, $false, $true, 1, $false, "compare the results of different", 2) | Out-Null
